$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-record data which gets permuted across rows 2-7
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Capture the current ("before") values for the affected rows/columns
$data = @{}
foreach ($r in 2..7) {
    $data[$r] = @{}
    foreach ($col in $cols) {
        $data[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# New row r receives the data that used to live in row $map[r]
$map = @{ 2 = 3; 3 = 2; 4 = 6; 5 = 7; 6 = 5; 7 = 4 }

foreach ($r in 2..7) {
    $src = $map[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $data[$src][$col]
    }
}
